# Update "想去人数" (F column) figures across the worksheets, per the
# upstream data refresh recorded in the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1576
$ws1.Range("F5").Value  = 273
$ws1.Range("F6").Value  = 60
$ws1.Range("F7").Value  = 1600
$ws1.Range("F8").Value  = 10213
$ws1.Range("F11").Value = 258
$ws1.Range("F12").Value = 194
$ws1.Range("F13").Value = 392
$ws1.Range("F14").Value = 7077
$ws1.Range("F15").Value = 1103
$ws1.Range("F17").Value = 37

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 9

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1576
$ws4.Range("F5").Value  = 273
$ws4.Range("F6").Value  = 9
$ws4.Range("F7").Value  = 60
$ws4.Range("F8").Value  = 1600
$ws4.Range("F11").Value = 10213
$ws4.Range("F14").Value = 258
$ws4.Range("F15").Value = 194
$ws4.Range("F16").Value = 392
$ws4.Range("F17").Value = 7077
$ws4.Range("F18").Value = 1103
$ws4.Range("F20").Value = 37
